$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$spaces77 = "".PadRight(77)

# ---------------------------------------------------------------------------
# Block 1: "To :" cell - insert Address: / blank / Contact Number: / blank
# paragraphs right after the empty paragraph that follows "To :" (before the
# Heading2 paragraph that holds the <w:br/>).
# ---------------------------------------------------------------------------
$anchor1 = $d.Paragraphs.Item(7)
$anchor1.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item(8)
$block1 = "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t>Address:</w:t></w:r></w:p>" +
          "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr></w:p>" +
          "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t>Contact Number:</w:t></w:r><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t xml:space='preserve'>$spaces77</w:t></w:r></w:p>" +
          "<w:p $wns><w:pPr><w:ind w:left='0'/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>"
$newPara1.Range.InsertXML($block1)

# ---------------------------------------------------------------------------
# Block 2: "From:" cell - trim the trailing spaces in the From: paragraph to
# a single space and give the paragraph mark the sz/szCs run properties,
# then insert blank / Address: / blank / Contact Number: paragraphs after it.
# After the 4 paragraphs inserted above, the From: paragraph shifted from
# index 13 to index 17.
# ---------------------------------------------------------------------------
$pFrom = $d.Paragraphs.Item(17)
$newXmlFrom = "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr><w:r><w:t xml:space='preserve'>           </w:t></w:r><w:r><w:t xml:space='preserve'>  </w:t></w:r><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t>From:</w:t></w:r><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r></w:p>"
$pFrom.Range.InsertXML($newXmlFrom)

$pFrom2 = $d.Paragraphs.Item(17)
$pFrom2.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(18)
$block2 = "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr></w:p>" +
          "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t xml:space='preserve'>          Address:</w:t></w:r></w:p>" +
          "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr></w:pPr></w:p>" +
          "<w:p $wns><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t xml:space='preserve'>          Contact Number:</w:t></w:r><w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='44'/></w:rPr><w:t xml:space='preserve'>$spaces77</w:t></w:r></w:p>"
$newPara2.Range.InsertXML($block2)

Write-Host "Edit complete"
